# DataSource - Emision Motor Flota Parte 1.xlsx
# Commit: "se arreglo rutas a C:temp se modifico el valor por defecto de la Franquicia Fija"
#
# Updates the sample/default row of data used by the Ranorex test data source:
#  - DatosMotor (sheet1): environment/URL/credentials + vehicle + policy defaults
#  - DatosPAS  (sheet2): group code + group name

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("DatosMotor")
$ws2 = $wb.Worksheets.Item("DatosPAS")

# --- DatosMotor!row 2 -------------------------------------------------
# Order matters: new shared strings get appended in the order they are
# first assigned, so write them in the same order the target file uses.

$ws1.Range("Q2").Value = "23/06/2021"
$ws1.Range("S2").Value = "TOYOTA"
$ws1.Range("T2").Value = "ETIOS 1.5 4 PTAS X 6MT L/18"

$ws2.Range("B2").Value = "Agustin Seisdedos"

$ws1.Range("A2").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws1.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws1.Range("D2").Value = "gw"
$ws1.Range("W2").Value = "JUN330"
$ws1.Range("X2").Value = "A1234567JUN330"
$ws1.Range("Y2").Value = "B1234567JUN330"

$ws1.Range("E2").Value = 5944085871
$ws1.Range("R2").Value = 2020
$ws1.Range("U2").Value = 1481000
$ws1.Range("V2").Value = "TR - Todo Riesgo Franquicia Fija"

# Re-establish the hyperlink on the URL cell (ambiente login URL)
[void]$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do")
# Hyperlinks.Add() clones the cell style into a new (duplicate) style
# record; re-apply the named style so B2 keeps referencing the original
# "Hipervinculo" style entry instead of a near-duplicate one.
$ws1.Range("B2").Style = "Hipervínculo"

# --- DatosPAS!row 2 -----------------------------------------------------
$ws2.Range("A2").Value = 6254

# --- Selections / active sheet ------------------------------------------
[void]$ws2.Activate()
[void]$ws2.Range("F11").Select()

[void]$ws1.Activate()
[void]$ws1.Range("Y7").Select()
